$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (46075 -> 46076) for every data row (rows 2 through 80).
$ws.Range("C2:C80").Value = 46076
